# test_reallocate.xlsx - "Accounts" sheet update
#
# Adds row/column SUM() totals (T2:T14 across the row, and E15:T15 down each
# column) and rounds a handful of "Cash/MMKT" balances (column E, rows 4-10)
# down to whole dollars - matching the new differenceCurrentDesiredAccounts /
# portfolioValue / categoryTotal helpers + their unit tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accounts")

# --- round a few Cash/MMKT (column E) balances down to whole dollars -------
$ws.Range("E4").Value2  = 13793
$ws.Range("E5").Value2  = 447
$ws.Range("E6").Value2  = 6
$ws.Range("E7").Value2  = 79
$ws.Range("E8").Value2  = 220
$ws.Range("E9").Value2  = 283
$ws.Range("E10").Value2 = 54

# --- row totals: T2:T14 = SUM(E:S) for that row -----------------------------
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("T$r").Formula = "=SUM(E" + $r + ":S" + $r + ")"
}

# --- grand-total row 15: SUM each column from row 2 through row 14 ---------
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
foreach ($col in $cols) {
    $ws.Range($col + "15").Formula = "=SUM(" + $col + "2:" + $col + "14)"
}

# match the formatting of the grand-total row (style carried from E15, which
# already used the medium-right-border look) across F15:T15
$ws.Range("F15:T15").Style = $ws.Range("E15").Style

# --- restore the saved selection -------------------------------------------
$ws.Range("E7").Select()
